$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 880.3939
$ws.Range("I33").Value = 676.2963
$ws.Range("J33").Value = 1798.8334
$ws.Range("K33").Value = 676.2963
$ws.Range("L33").Value = 1798.8334
$ws.Range("M33").Value = -447.2963
$ws.Range("N33").Value = -2256.8334

$ws.Range("H40").Value = 3411.9644
$ws.Range("I40").Value = 2216.875
$ws.Range("J40").Value = 3890
$ws.Range("K40").Value = 2216.875
$ws.Range("L40").Value = 3890
$ws.Range("M40").Value = -2041.875
$ws.Range("N40").Value = -4240

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 8340.825000000001
$ws.Range("I132").Value = 5983.0225
$ws.Range("J132").Value = 14235.333
$ws.Range("K132").Value = 17949.0675
$ws.Range("L132").Value = 42705.999
$ws.Range("M132").Value = -15419.0675
$ws.Range("N132").Value = -47765.999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 5066.3335
$ws.Range("I5").Value = 15000
$ws.Range("J5").Value = 99.5
$ws.Range("K5").Value = 15000
$ws.Range("L5").Value = 99.5
$ws.Range("M5").Value = -14887
$ws.Range("N5").Value = -325.5

$ws.Range("H102").Value = 70368.28999999999
$ws.Range("I102").Value = 30999.666
$ws.Range("J102").Value = 99894.75
$ws.Range("K102").Value = 30999.666
$ws.Range("L102").Value = 99894.75
$ws.Range("M102").Value = -27754.666
$ws.Range("N102").Value = -106384.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 293.42856
$ws.Range("I22").Value = 261
$ws.Range("J22").Value = 488
$ws.Range("K22").Value = 261
$ws.Range("L22").Value = 488
$ws.Range("M22").Value = 89
$ws.Range("N22").Value = -1188

$ws.Range("H31").Value = 734450.6
$ws.Range("I31").Value = 5450.8
$ws.Range("J31").Value = 1065814.2
$ws.Range("K31").Value = 5450.8
$ws.Range("L31").Value = 1065814.2
$ws.Range("M31").Value = -5155.8
$ws.Range("N31").Value = -1066404.2

$ws.Range("H34").Value = 734450.6
$ws.Range("I34").Value = 5450.8
$ws.Range("J34").Value = 1065814.2
$ws.Range("K34").Value = 5450.8
$ws.Range("L34").Value = 1065814.2
$ws.Range("M34").Value = -5248.8
$ws.Range("N34").Value = -1066218.2

$ws.Range("H41").Value = 5583
$ws.Range("I41").Value = 5583
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 5583
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -5155

$ws.Range("H51").Value = 46599.9
$ws.Range("I51").Value = 36571.285
$ws.Range("J51").Value = 70000
$ws.Range("K51").Value = 36571.285
$ws.Range("L51").Value = 70000
$ws.Range("M51").Value = -35835.285
$ws.Range("N51").Value = -71472

$ws.Range("H60").Value = 74450
$ws.Range("I60").Value = 77250
$ws.Range("J60").Value = 73750
$ws.Range("K60").Value = 77250
$ws.Range("L60").Value = 73750
$ws.Range("M60").Value = -76739
$ws.Range("N60").Value = -74772

$ws.Range("H61").Value = 46599.9
$ws.Range("I61").Value = 36571.285
$ws.Range("J61").Value = 70000
$ws.Range("K61").Value = 36571.285
$ws.Range("L61").Value = 70000
$ws.Range("M61").Value = -36223.285
$ws.Range("N61").Value = -70696

$ws.Range("H68").Value = 74980
$ws.Range("I68").Value = 74980
$ws.Range("J68").Value = 74980
$ws.Range("K68").Value = 74980
$ws.Range("L68").Value = 74980
$ws.Range("M68").Value = -74231
$ws.Range("N68").Value = -76478

$ws.Range("H71").Value = 74980
$ws.Range("I71").Value = 74980
$ws.Range("J71").Value = 74980
$ws.Range("K71").Value = 224940
$ws.Range("L71").Value = 224940
$ws.Range("M71").Value = -221196
$ws.Range("N71").Value = -232428

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2284.5
$ws.Range("I5").Value = 380
$ws.Range("J5").Value = 7998
$ws.Range("K5").Value = 1140
$ws.Range("L5").Value = 23994
$ws.Range("M5").Value = -1028
$ws.Range("N5").Value = -24218

$ws.Range("H62").Value = 14000
$ws.Range("I62").Value = 14000
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 42000
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -41314

$ws.Range("H65").Value = 14000
$ws.Range("I65").Value = 14000
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 126000
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -122568

$ws.Range("H82").Value = 8485.571
$ws.Range("I82").Value = 4949.5
$ws.Range("J82").Value = 9900
$ws.Range("K82").Value = 14848.5
$ws.Range("L82").Value = 29700
$ws.Range("M82").Value = -14442.5
$ws.Range("N82").Value = -30512

$ws.Range("H85").Value = 8485.571
$ws.Range("I85").Value = 4949.5
$ws.Range("J85").Value = 9900
$ws.Range("K85").Value = 14848.5
$ws.Range("L85").Value = 29700
$ws.Range("M85").Value = -13444.5
$ws.Range("N85").Value = -32508

$ws.Range("H95").Value = 10000
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 10000
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 30000
$ws.Range("M95").ClearContents()
$ws.Range("N95").Value = -34118

$ws.Range("H122").Value = 2344.1667
$ws.Range("I122").Value = 1200
$ws.Range("J122").Value = 2573
$ws.Range("K122").Value = 10800
$ws.Range("L122").Value = 23157
$ws.Range("M122").Value = -8350
$ws.Range("N122").Value = -28057

$ws.Range("H132").Value = 824.7778
$ws.Range("I132").Value = 824.7778
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7423.000199999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4893.000199999999
$ws.Range("N132").ClearContents()

$ws.Range("H134").Value = 3484.1428
$ws.Range("I134").Value = 2355.3215
$ws.Range("J134").Value = 7999.4287
$ws.Range("K134").Value = 7065.9645
$ws.Range("L134").Value = 23998.2861
$ws.Range("M134").Value = -1995.9645
$ws.Range("N134").Value = -34138.2861

$ws.Range("H135").Value = 2284.5
$ws.Range("I135").Value = 380
$ws.Range("J135").Value = 7998
$ws.Range("K135").Value = 3420
$ws.Range("L135").Value = 71982
$ws.Range("M135").Value = -885
$ws.Range("N135").Value = -77052

$ws.Range("H140").Value = 160550.73
$ws.Range("I140").Value = 169220.22
$ws.Range("J140").Value = 4500
$ws.Range("K140").Value = 507660.66
$ws.Range("L140").Value = 13500
$ws.Range("M140").Value = -502480.66
$ws.Range("N140").Value = -23860

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 608.8333
$ws.Range("I3").Value = 251
$ws.Range("J3").Value = 966.6667
$ws.Range("K3").Value = 251
$ws.Range("L3").Value = 966.6667
$ws.Range("M3").Value = -135
$ws.Range("N3").Value = -1198.6667

$ws.Range("H7").Value = 61001.8
$ws.Range("I7").Value = 5000
$ws.Range("J7").Value = 75002.25
$ws.Range("K7").Value = 5000
$ws.Range("L7").Value = 75002.25
$ws.Range("M7").Value = -4888
$ws.Range("N7").Value = -75226.25

$ws.Range("H8").Value = 61001.8
$ws.Range("I8").Value = 5000
$ws.Range("J8").Value = 75002.25
$ws.Range("K8").Value = 5000
$ws.Range("L8").Value = 75002.25
$ws.Range("M8").Value = -4861
$ws.Range("N8").Value = -75280.25

$ws.Range("H10").Value = 3450
$ws.Range("I10").Value = 3450
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 3450
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -3281
$ws.Range("N10").ClearContents()

$ws.Range("H11").Value = 15040215
$ws.Range("I11").Value = 8257111
$ws.Range("J11").Value = 27249802
$ws.Range("K11").Value = 8257111
$ws.Range("L11").Value = 27249802
$ws.Range("M11").Value = -8256972
$ws.Range("N11").Value = -27250080

$ws.Range("H12").Value = 15000
$ws.Range("I12").Value = 10000
$ws.Range("J12").Value = 20000
$ws.Range("K12").Value = 10000
$ws.Range("L12").Value = 20000
$ws.Range("M12").Value = -9860
$ws.Range("N12").Value = -20280

$ws.Range("H109").Value = 47749
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 47749
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 47749
$ws.Range("N109").Value = -49829

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

$ws.Range("H9").Value = 566.6667
$ws.Range("I9").Value = 566.6667
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 566.6667
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -342.6667

$ws.Range("H21").Value = 9728.444
$ws.Range("I21").Value = 10000
$ws.Range("J21").Value = 9712.471
$ws.Range("K21").Value = 10000
$ws.Range("L21").Value = 9712.471
$ws.Range("M21").Value = -9826
$ws.Range("N21").Value = -10060.471

$ws.Range("H23").Value = 27450
$ws.Range("I23").Value = 27450
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 27450
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -27220

$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()

$ws.Range("H34").Value = 7499.5
$ws.Range("I34").Value = 7499.5
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 7499.5
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -7327.5

$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()

$ws.Range("H43").Value = 1849809.1
$ws.Range("I43").Value = 2249766.8
$ws.Range("J43").Value = 50000
$ws.Range("K43").Value = 2249766.8
$ws.Range("L43").Value = 50000
$ws.Range("M43").Value = -2249573.8
$ws.Range("N43").Value = -50386

$ws.Range("H46").Value = 3460.1428
$ws.Range("I46").Value = 3361.4375
$ws.Range("J46").Value = 3591.75
$ws.Range("K46").Value = 3361.4375
$ws.Range("L46").Value = 3591.75
$ws.Range("M46").Value = -3173.4375
$ws.Range("N46").Value = -3967.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 5502.5
$ws.Range("I6").Value = 999
$ws.Range("J6").Value = 10006
$ws.Range("K6").Value = 999
$ws.Range("L6").Value = 10006
$ws.Range("M6").Value = -884
$ws.Range("N6").Value = -10236

$ws.Range("H23").Value = 688
$ws.Range("I23").Value = 688
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 688
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -459
$ws.Range("N23").ClearContents()

$ws.Range("H41").Value = 14825.667
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 14825.667
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 14825.667
$ws.Range("N41").Value = -15605.667
